$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2").Value = "In Translation"
$wsOverview.Range("F2").Value = "In Translation"

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C2").Value = "In Translation"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C2").Value = "In Translation"

$wsOverview.Range("E:F").Columns.AutoFit() | Out-Null
$wsZhCn.Range("C:C").Columns.AutoFit() | Out-Null
$wsDeDe.Range("C:C").Columns.AutoFit() | Out-Null
